$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GSM")
$ws.Columns("D").Insert()
$ws.Range("D7").Value = 43373
